$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '42.015.53'
$ws.Range('E2').Value = '  +5.20%  '
$ws.Range('D3').Value = '2.257.55'
$ws.Range('E3').Value = '  +1.79%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '301.74'
$ws.Range('E5').Value = '  +3.39%  '
$ws.Range('D6').Value = '92.84'
$ws.Range('E6').Value = '  +5.88%  '
$ws.Range('D7').Value = '0.532'
$ws.Range('E7').Value = '  +3.63%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +3.08%  '
$ws.Range('D10').Value = '32.73'
$ws.Range('E10').Value = '  +7.12%  '
$ws.Range('D11').Value = '54.55'
$ws.Range('E11').Value = '  +8.97%  '
$ws.Range('E12').Value = '  +2.47%  '
$ws.Range('E13').Value = '  +3.36%  '
$ws.Range('E14').Value = '  +3.60%  '
$ws.Range('D15').Value = '2.602.08'
$ws.Range('E15').Value = '  +1.65%  '
$ws.Range('D16').Value = '14.16'
$ws.Range('E16').Value = '  +2.78%  '
$ws.Range('D17').Value = '2.259.00'
$ws.Range('E17').Value = '  +3.89%  '
$ws.Range('E18').Value = '  +3.50%  '
$ws.Range('D19').Value = '41.889.57'
$ws.Range('E19').Value = '  +5.06%  '
$ws.Range('D20').Value = '12.20'
$ws.Range('E20').Value = '  +9.72%  '
$ws.Range('D21').Value = '0.0₃0907'
$ws.Range('E21').Value = '  +2.20%  '
$ws.Range('E22').Value = '  +3.74%  '
$ws.Range('E23').Value = '  +2.41%  '
$ws.Range('D24').Value = '241.91'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('E25').Value = '  +5.19%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +4.15%  '
$ws.Range('D28').Value = '23.97'
$ws.Range('E28').Value = '  +3.40%  '
$ws.Range('D29').Value = '9.70'
$ws.Range('E29').Value = '  +4.88%  '
$ws.Range('D31').Value = '34.12'
$ws.Range('E31').Value = '  +6.68%  '
$ws.Range('D32').Value = '158.98'
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('D33').Value = '0.999'
$ws.Range('D34').Value = '5.16'
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').Value = '0.0743'
$ws.Range('E35').Value = '  +4.32%  '
$ws.Range('D36').Value = '3.06'
$ws.Range('E36').Value = '  +3.22%  '
$ws.Range('E37').Value = '  +2.51%  '
$ws.Range('E38').Value = '  +5.61%  '
$ws.Range('D39').Value = '16.60'
$ws.Range('E39').Value = '  +8.31%  '
$ws.Range('E40').Value = '  +3.97%  '
$ws.Range('E41').Value = '  +4.29%  '
$ws.Range('E42').Value = '  +5.72%  '
$ws.Range('D43').Value = '2.052.65'
$ws.Range('E43').Value = '  -2.63%  '
$ws.Range('D44').Value = '19.95'
$ws.Range('E44').Value = '  +11.76%  '
$ws.Range('D45').Value = '0.0280'
$ws.Range('E45').Value = '  +3.72%  '
$ws.Range('E46').Value = '  +1.98%  '
$ws.Range('E47').Value = '  +7.73%  '
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('D49').Value = '2.475.55'
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('E50').Value = '  +2.88%  '
$ws.Range('E51').Value = '  +4.48%  '
